$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = 'Datos actualizados a 25 de Marzo de 2020 a las 14:16'

# Refresh the country data: the source feed re-ranked a number of
# countries (by total cases) and updated several counters.

$ws.Cells.Item(6, 1).Value = 'Estados Unidos'
$ws.Cells.Item(6, 2).Value = 54999
$ws.Cells.Item(6, 3).Value = 143
$ws.Cells.Item(6, 4).Value = 379
$ws.Cells.Item(6, 5).Value = 53835
$ws.Cells.Item(6, 6).Value = 1175
$ws.Cells.Item(6, 7).Value = 5
$ws.Cells.Item(6, 8).Value = 785

$ws.Cells.Item(8, 1).Value = 'Alemania'
$ws.Cells.Item(8, 2).Value = 34055
$ws.Cells.Item(8, 3).Value = 1064
$ws.Cells.Item(8, 4).Value = 3540
$ws.Cells.Item(8, 5).Value = 30340
$ws.Cells.Item(8, 6).Value = 23
$ws.Cells.Item(8, 7).Value = 16
$ws.Cells.Item(8, 8).Value = 175

$ws.Cells.Item(15, 1).Value = 'Austria'
$ws.Cells.Item(15, 2).Value = 5516
$ws.Cells.Item(15, 3).Value = 233
$ws.Cells.Item(15, 4).Value = 9
$ws.Cells.Item(15, 5).Value = 5477
$ws.Cells.Item(15, 6).Value = 26
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(15, 8).Value = 30

$ws.Cells.Item(18, 1).Value = 'Noruega'
$ws.Cells.Item(18, 2).Value = 2902
$ws.Cells.Item(18, 3).Value = 36
$ws.Cells.Item(18, 4).Value = 6
$ws.Cells.Item(18, 5).Value = 2882
$ws.Cells.Item(18, 6).Value = 57
$ws.Cells.Item(18, 7).Value = 2
$ws.Cells.Item(18, 8).Value = 14

$ws.Cells.Item(32, 1).Value = 'Pakistan'
$ws.Cells.Item(32, 2).Value = 1014
$ws.Cells.Item(32, 3).Value = 42
$ws.Cells.Item(32, 4).Value = 21
$ws.Cells.Item(32, 5).Value = 985
$ws.Cells.Item(32, 6).Value = 5
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(32, 8).Value = 8

$ws.Cells.Item(35, 1).Value = 'Chile'
$ws.Cells.Item(35, 2).Value = 922
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 22
$ws.Cells.Item(35, 5).Value = 898
$ws.Cells.Item(35, 6).Value = 7
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 2

$ws.Cells.Item(37, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(37, 2).Value = 900
$ws.Cells.Item(37, 3).Value = 133
$ws.Cells.Item(37, 4).Value = 29
$ws.Cells.Item(37, 5).Value = 869
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 1
$ws.Cells.Item(37, 8).Value = 2

$ws.Cells.Item(51, 1).Value = 'Egipto'
$ws.Cells.Item(51, 2).Value = 442
$ws.Cells.Item(51, 3).Value = 40
$ws.Cells.Item(51, 4).Value = 80
$ws.Cells.Item(51, 5).Value = 342
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 20

$ws.Cells.Item(52, 1).Value = 'Barein'
$ws.Cells.Item(52, 2).Value = 419
$ws.Cells.Item(52, 3).Value = 27
$ws.Cells.Item(52, 4).Value = 177
$ws.Cells.Item(52, 5).Value = 239
$ws.Cells.Item(52, 6).Value = 2
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 3

$ws.Cells.Item(53, 1).Value = 'Croacia'
$ws.Cells.Item(53, 2).Value = 418
$ws.Cells.Item(53, 3).Value = 36
$ws.Cells.Item(53, 4).Value = 16
$ws.Cells.Item(53, 5).Value = 401
$ws.Cells.Item(53, 6).Value = 6
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 1

$ws.Cells.Item(54, 1).Value = 'Peru'
$ws.Cells.Item(54, 2).Value = 416
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 1
$ws.Cells.Item(54, 5).Value = 408
$ws.Cells.Item(54, 6).Value = 9
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 7

$ws.Cells.Item(55, 1).Value = 'Hong Kong'
$ws.Cells.Item(55, 2).Value = 410
$ws.Cells.Item(55, 3).Value = 23
$ws.Cells.Item(55, 4).Value = 102
$ws.Cells.Item(55, 5).Value = 304
$ws.Cells.Item(55, 6).Value = 4
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 4

$ws.Cells.Item(56, 1).Value = 'Mexico'
$ws.Cells.Item(56, 2).Value = 405
$ws.Cells.Item(56, 3).Value = 38
$ws.Cells.Item(56, 4).Value = 4
$ws.Cells.Item(56, 5).Value = 396
$ws.Cells.Item(56, 6).Value = 1
$ws.Cells.Item(56, 7).Value = 1
$ws.Cells.Item(56, 8).Value = 5

$ws.Cells.Item(57, 1).Value = 'Estonia'
$ws.Cells.Item(57, 2).Value = 404
$ws.Cells.Item(57, 3).Value = 35
$ws.Cells.Item(57, 4).Value = 8
$ws.Cells.Item(57, 5).Value = 396
$ws.Cells.Item(57, 6).Value = 5
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 0

$ws.Cells.Item(61, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(61, 2).Value = 333
$ws.Cells.Item(61, 3).Value = 85
$ws.Cells.Item(61, 4).Value = 45
$ws.Cells.Item(61, 5).Value = 286
$ws.Cells.Item(61, 6).Value = 2
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 2

$ws.Cells.Item(62, 1).Value = 'Irak'
$ws.Cells.Item(62, 2).Value = 316
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 75
$ws.Cells.Item(62, 5).Value = 214
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 27

$ws.Cells.Item(63, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(63, 2).Value = 312
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 3
$ws.Cells.Item(63, 5).Value = 303
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 6

$ws.Cells.Item(64, 1).Value = 'Serbia'
$ws.Cells.Item(64, 2).Value = 303
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(64, 4).Value = 15
$ws.Cells.Item(64, 5).Value = 284
$ws.Cells.Item(64, 6).Value = 21
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 4

$ws.Cells.Item(65, 1).Value = 'Armenia'
$ws.Cells.Item(65, 2).Value = 265
$ws.Cells.Item(65, 3).Value = 16
$ws.Cells.Item(65, 4).Value = 16
$ws.Cells.Item(65, 5).Value = 249
$ws.Cells.Item(65, 6).Value = 6
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0

$ws.Cells.Item(66, 1).Value = 'Argelia'
$ws.Cells.Item(66, 2).Value = 264
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = 65
$ws.Cells.Item(66, 5).Value = 180
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 19

$ws.Cells.Item(67, 1).Value = 'Lituania'
$ws.Cells.Item(67, 2).Value = 255
$ws.Cells.Item(67, 3).Value = 46
$ws.Cells.Item(67, 4).Value = 1
$ws.Cells.Item(67, 5).Value = 250
$ws.Cells.Item(67, 6).Value = 1
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = 4

$ws.Cells.Item(76, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(76, 2).Value = 188
$ws.Cells.Item(76, 3).Value = 24
$ws.Cells.Item(76, 4).Value = 1
$ws.Cells.Item(76, 5).Value = 186
$ws.Cells.Item(76, 6).Value = 6
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 1

$ws.Cells.Item(77, 1).Value = 'San Marino'
$ws.Cells.Item(77, 2).Value = 187
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 4
$ws.Cells.Item(77, 5).Value = 162
$ws.Cells.Item(77, 6).Value = 12
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 21

$ws.Cells.Item(78, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(78, 2).Value = 177
$ws.Cells.Item(78, 3).Value = 29
$ws.Cells.Item(78, 4).Value = 1
$ws.Cells.Item(78, 5).Value = 174
$ws.Cells.Item(78, 6).Value = 1
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 2

$ws.Cells.Item(79, 1).Value = 'Costa Rica'
$ws.Cells.Item(79, 2).Value = 177
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(79, 4).Value = 2
$ws.Cells.Item(79, 5).Value = 173
$ws.Cells.Item(79, 6).Value = 4
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 2

$ws.Cells.Item(80, 1).Value = 'Marruecos'
$ws.Cells.Item(80, 2).Value = 170
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(80, 4).Value = 6
$ws.Cells.Item(80, 5).Value = 159
$ws.Cells.Item(80, 6).Value = 1
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 5

$ws.Cells.Item(81, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(81, 2).Value = 168
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = 2
$ws.Cells.Item(81, 5).Value = 163
$ws.Cells.Item(81, 6).Value = 1
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 3

$ws.Cells.Item(97, 1).Value = 'Camboya'
$ws.Cells.Item(97, 2).Value = 93
$ws.Cells.Item(97, 3).Value = 2
$ws.Cells.Item(97, 4).Value = 6
$ws.Cells.Item(97, 5).Value = 87
$ws.Cells.Item(97, 6).Value = 1
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0

$ws.Cells.Item(107, 1).Value = 'Ghana'
$ws.Cells.Item(107, 2).Value = 68
$ws.Cells.Item(107, 3).Value = 15
$ws.Cells.Item(107, 4).Value = 0
$ws.Cells.Item(107, 5).Value = 65
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 3

$ws.Cells.Item(114, 1).Value = 'Cuba'
$ws.Cells.Item(114, 2).Value = 48
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 1
$ws.Cells.Item(114, 5).Value = 46
$ws.Cells.Item(114, 6).Value = 2
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 1

$ws.Cells.Item(115, 1).Value = 'Mauricio'
$ws.Cells.Item(115, 2).Value = 48
$ws.Cells.Item(115, 3).Value = 6
$ws.Cells.Item(115, 4).Value = 0
$ws.Cells.Item(115, 5).Value = 46
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 2

$ws.Cells.Item(116, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(116, 2).Value = 48
$ws.Cells.Item(116, 3).Value = 3
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 5).Value = 46
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 2

$ws.Cells.Item(123, 1).Value = 'Mayotte'
$ws.Cells.Item(123, 2).Value = 36
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 36
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 0

$ws.Cells.Item(124, 1).Value = 'Honduras'
$ws.Cells.Item(124, 2).Value = 36
$ws.Cells.Item(124, 3).Value = 6
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(124, 5).Value = 36
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 0

$ws.Cells.Item(141, 1).Value = 'Uganda'
$ws.Cells.Item(141, 2).Value = 14
$ws.Cells.Item(141, 3).Value = 5
$ws.Cells.Item(141, 4).Value = 0
$ws.Cells.Item(141, 5).Value = 14
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 0

$ws.Cells.Item(142, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(142, 2).Value = 14
$ws.Cells.Item(142, 3).Value = 4
$ws.Cells.Item(142, 4).Value = 0
$ws.Cells.Item(142, 5).Value = 14
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 0

$ws.Cells.Item(144, 1).Value = 'Etiopia'
$ws.Cells.Item(144, 2).Value = 12
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 12
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 0

$ws.Cells.Item(146, 1).Value = 'Tanzania'
$ws.Cells.Item(146, 2).Value = 12
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 0
$ws.Cells.Item(146, 5).Value = 12
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 0

$ws.Cells.Item(153, 1).Value = 'Surinam'
$ws.Cells.Item(153, 2).Value = 7
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 7
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 0

$ws.Cells.Item(154, 1).Value = 'Haiti'
$ws.Cells.Item(154, 2).Value = 7
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 5).Value = 7
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 0

$ws.Cells.Item(155, 1).Value = 'Dominica'
$ws.Cells.Item(155, 2).Value = 7
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 0
$ws.Cells.Item(155, 5).Value = 7
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 0

$ws.Cells.Item(158, 1).Value = 'Bermudas'
$ws.Cells.Item(158, 2).Value = 6
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 6
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 0

$ws.Cells.Item(159, 1).Value = 'Benin'
$ws.Cells.Item(159, 2).Value = 6
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 6
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

$ws.Cells.Item(160, 1).Value = 'Curazao'
$ws.Cells.Item(160, 2).Value = 6
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 0
$ws.Cells.Item(160, 5).Value = 5
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 1

$ws.Cells.Item(162, 1).Value = 'Gabon'
$ws.Cells.Item(162, 2).Value = 6
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 0
$ws.Cells.Item(162, 5).Value = 5
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 1

$ws.Cells.Item(167, 1).Value = 'Congo'
$ws.Cells.Item(167, 2).Value = 4
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 0
$ws.Cells.Item(167, 5).Value = 4
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 0

$ws.Cells.Item(168, 1).Value = 'Suazilandia'
$ws.Cells.Item(168, 2).Value = 4
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 5).Value = 4
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 0

$ws.Cells.Item(172, 1).Value = 'Angola'
$ws.Cells.Item(172, 2).Value = 3
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 3
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(173, 1).Value = 'Santa Lucia'
$ws.Cells.Item(173, 2).Value = 3
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 3
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(174, 1).Value = 'Liberia'
$ws.Cells.Item(174, 2).Value = 3
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 3
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

$ws.Cells.Item(175, 1).Value = 'San Bartolome'
$ws.Cells.Item(175, 2).Value = 3
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 3
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

$ws.Cells.Item(176, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(176, 2).Value = 3
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 3
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 1).Value = 'Republica del Chad'
$ws.Cells.Item(177, 2).Value = 3
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 3
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

$ws.Cells.Item(178, 1).Value = 'Birmania'
$ws.Cells.Item(178, 2).Value = 3
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 3
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 0

$ws.Cells.Item(179, 1).Value = 'Mozambique'
$ws.Cells.Item(179, 2).Value = 3
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 3
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0

$ws.Cells.Item(180, 1).Value = 'Laos'
$ws.Cells.Item(180, 2).Value = 3
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 3
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(181, 2).Value = 3
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 3
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

$ws.Cells.Item(182, 1).Value = 'Gambia'
$ws.Cells.Item(182, 2).Value = 3
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 2
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 1

$ws.Cells.Item(183, 1).Value = 'Nepal'
$ws.Cells.Item(183, 2).Value = 3
$ws.Cells.Item(183, 3).Value = 1
$ws.Cells.Item(183, 4).Value = 1
$ws.Cells.Item(183, 5).Value = 2
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

$ws.Cells.Item(184, 1).Value = 'Zimbabue'
$ws.Cells.Item(184, 2).Value = 3
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 2
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 1

$ws.Cells.Item(187, 1).Value = 'Mauritania'
$ws.Cells.Item(187, 2).Value = 2
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 2
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

$ws.Cells.Item(188, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(188, 2).Value = 2
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 2
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

$ws.Cells.Item(189, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(189, 2).Value = 2
$ws.Cells.Item(189, 3).Value = 2
$ws.Cells.Item(189, 4).Value = 0
$ws.Cells.Item(189, 5).Value = 2
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 0

$ws.Cells.Item(190, 1).Value = 'Mali'
$ws.Cells.Item(190, 2).Value = 2
$ws.Cells.Item(190, 3).Value = 2
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 2
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

$ws.Cells.Item(191, 1).Value = 'Nicaragua'
$ws.Cells.Item(191, 2).Value = 2
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 2
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 0

$ws.Cells.Item(192, 1).Value = 'Montserrat'
$ws.Cells.Item(192, 2).Value = 1
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 1
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

$ws.Cells.Item(193, 1).Value = 'Siria'
$ws.Cells.Item(193, 2).Value = 1
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 5).Value = 1
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(194, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(194, 2).Value = 1
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 5).Value = 1
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

$ws.Cells.Item(195, 1).Value = 'Belice'
$ws.Cells.Item(195, 2).Value = 1
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 1
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

$ws.Cells.Item(196, 1).Value = 'Eritrea'
$ws.Cells.Item(196, 2).Value = 1
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 0
$ws.Cells.Item(196, 5).Value = 1
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 0

$ws.Cells.Item(197, 1).Value = 'Timor Oriental'
$ws.Cells.Item(197, 2).Value = 1
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 5).Value = 1
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0

$ws.Cells.Item(198, 1).Value = 'Libia'
$ws.Cells.Item(198, 2).Value = 1
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 5).Value = 1
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0

$ws.Cells.Item(199, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(199, 2).Value = 1
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 0
$ws.Cells.Item(199, 5).Value = 1
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0
